# Removes the "raw xml" "@" marker from the {@Arsak} placeholder so it
# reads {Arsak}, and relocates the (hidden) _GoBack bookmark to the edit
# point, matching native Word's "last edit location" bookmark behaviour.

$d = $word.ActiveDocument

# 1) Drop the old _GoBack bookmark first so the name is free to reuse.
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# 2) Locate the literal "@" run inside the "{@Arsak}" placeholder.
$rng = $d.Content
$found = $rng.Find.Execute("@")
if (-not $found) {
    throw "Could not find the '@' placeholder marker"
}

# 3) Remove the "@" character.
$rng.Delete()

# 4) Re-create _GoBack, collapsed, at the position the edit just happened
#    (exactly where Word itself would drop it after an in-place deletion).
$rng.Collapse(1)  # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
